# Update the LR-pairs data table (Efna2-Epha5) with new TPM-derived values.
# Target table is 10 data rows (A2:T11); columns:
#  A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#  E..T numeric metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efna2", "Epha5", "MuSCs", 3, 1, 3.293025333333333, 9.879076, 0.2347004583311417, 0.2347004583311417, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.08037177163555553, 0.7233459447199999, 0.2319945940411777, 0.2319945940411777),
    @("ECs", "Efna2", "Epha5", "Resolving-Mac", 3, 1, 3.293025333333333, 9.879076, 0.2347004583311417, 0.2347004583311417, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.0009374145448888888, 0.008436730904, 0.002705864289964023, 0.002705864289964023),
    @("FAPs", "Efna2", "Epha5", "MuSCs", 3, 1, 5.925454666666667, 17.776364, 0.422318927221656, 0.422318927221656, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.1446205968977778, 1.30158537208, 0.4174500074407977, 0.4174500074407978),
    @("FAPs", "Efna2", "Epha5", "Resolving-Mac", 3, 1, 5.925454666666667, 17.776364, 0.422318927221656, 0.422318927221656, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.001686779428444445, 0.015181014856, 0.004868919780858253, 0.004868919780858253),
    @("Inflammatory-Mac", "Efna2", "Epha5", "MuSCs", 3, 1, 2.186484333333333, 6.559453, 0.1558350827042511, 0.1558350827042511, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.05336479429555554, 0.4802831486599999, 0.1540384582391294, 0.1540384582391294),
    @("Inflammatory-Mac", "Efna2", "Epha5", "Resolving-Mac", 3, 1, 2.186484333333333, 6.559453, 0.1558350827042511, 0.1558350827042511, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.0006224192068888889, 0.005601772862, 0.001796624465121777, 0.001796624465121777),
    @("MuSCs", "Efna2", "Epha5", "MuSCs", 3, 1, 1.169741, 3.509223, 0.08336976519729013, 0.08336976519729013, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.02854947867333333, 0.25694530806, 0.08240859421315959, 0.08240859421315959),
    @("MuSCs", "Efna2", "Epha5", "Resolving-Mac", 3, 1, 1.169741, 3.509223, 0.08336976519729013, 0.08336976519729013, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.0003329862713333333, 0.002996876442, 0.0009611709841305422, 0.0009611709841305422),
    @("Resolving-Mac", "Efna2", "Epha5", "MuSCs", 3, 1, 1.456052666666666, 4.368157999999999, 0.1037757665456611, 0.1037757665456611, 3, 1, 0.02440666666666666, 0.07321999999999999, 0.9884709884709885, 0.9884709884709886, 0.03553739208444443, 0.3198365287599999, 0.1025793345367241, 0.1025793345367242),
    @("Resolving-Mac", "Efna2", "Epha5", "Resolving-Mac", 3, 1, 1.456052666666666, 4.368157999999999, 0.1037757665456611, 0.1037757665456611, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854, 0.01152901152901153, 0.01152901152901153, 0.000414489659111111, 0.003730406932, 0.001196432008936936, 0.001196432008936936)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
